# Update "想去人数" (want-to-go count) figures in the F column.
# Sheet 1 = "展览" (Exhibition), Sheet 4 = "全部类型" (All types) both carry
# the same exhibition rows (Sheet 4 has one extra row pulled in from the
# "演出" sheet at row 31, which is why its row numbers are offset by one
# starting there).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F5").Value = 983
$ws1.Range("F6").Value = 74
$ws1.Range("F7").Value = 2127
$ws1.Range("F10").Value = 4711
$ws1.Range("F13").Value = 295
$ws1.Range("F15").Value = 23
$ws1.Range("F16").Value = 160
$ws1.Range("F20").Value = 3627
$ws1.Range("F21").Value = 313
$ws1.Range("F22").Value = 587
$ws1.Range("F30").Value = 217
$ws1.Range("F31").Value = 18
$ws1.Range("F32").Value = 794
$ws1.Range("F33").Value = 2250
$ws1.Range("F34").Value = 414

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F5").Value = 983
$ws4.Range("F6").Value = 74
$ws4.Range("F7").Value = 2127
$ws4.Range("F10").Value = 4711
$ws4.Range("F13").Value = 295
$ws4.Range("F15").Value = 23
$ws4.Range("F16").Value = 160
$ws4.Range("F20").Value = 3627
$ws4.Range("F21").Value = 313
$ws4.Range("F22").Value = 587
$ws4.Range("F30").Value = 217
$ws4.Range("F32").Value = 18
$ws4.Range("F33").Value = 794
$ws4.Range("F34").Value = 2251
$ws4.Range("F35").Value = 414
